$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 8 and 9 (columns A, B, E, F, G, H, Q, R, Z, AB) ---
$ws.Range("A8").Value  = 131033281
$ws.Range("B8").Value  = 79862
$ws.Range("E8").Value  = 6453
$ws.Range("F8").Value  = "Vedskivlav"
$ws.Range("G8").Value  = "Hertelidea botryosa"
$ws.Range("H8").Value  = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q8").Value  = 395662
$ws.Range("R8").Value  = 6804783
$ws.Range("Z8").Value  = "11:57"
$ws.Range("AB8").Value = "11:57"

$ws.Range("A9").Value  = 131033360
$ws.Range("B9").Value  = 78909
$ws.Range("E9").Value  = 353
$ws.Range("F9").Value  = "Dvärgbägarlav"
$ws.Range("G9").Value  = "Cladonia parasitica"
$ws.Range("H9").Value  = "(Hoffm.) Hoffm."
$ws.Range("Q9").Value  = 395791
$ws.Range("R9").Value  = 6804722
$ws.Range("Z9").Value  = "13:19"
$ws.Range("AB9").Value = "13:19"

# --- Swap rows 17 and 18 (columns A, Q, R, Z, AB) ---
$ws.Range("A17").Value  = 131033336
$ws.Range("Q17").Value  = 395777
$ws.Range("R17").Value  = 6804741
$ws.Range("Z17").Value  = "13:19"
$ws.Range("AB17").Value = "13:19"

$ws.Range("A18").Value  = 131033318
$ws.Range("Q18").Value  = 395666
$ws.Range("R18").Value  = 6804680
$ws.Range("Z18").Value  = "12:00"
$ws.Range("AB18").Value = "12:00"
